# إضافة حدث جديد في Card12
# Row 20 was the last logged service event, stored with blank placeholder
# cells (B:K, N) instead of the usual "nan" text used elsewhere in the
# sheet. Adding the new event normalizes that previous row to "nan" and
# appends a fresh row 21 for the new event (date/servicer/correction),
# leaving its own unused columns blank - matching the sheet's convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card12")

# --- Normalize row 20's blank placeholder cells to "nan" ---
$blankCols20 = @(2,3,4,5,6,7,8,9,10,11,14)   # B,C,D,E,F,G,H,I,J,K,N
foreach ($col in $blankCols20) {
    $ws.Cells.Item(20, $col).Value = "nan"
}

# --- Append the new event as row 21 ---
# Columns that stay blank for this event (B..K and N), but still need to be
# *text* typed cells (matching every other cell on the sheet) rather than
# left as untouched/empty numeric cells. A leading quote forces Excel's
# text interpretation for an otherwise-empty value; the style is then reset
# to "Normal" so no stray quote-prefix number format lingers behind.
$blankCols21 = @(1,2,3,4,5,6,7,8,9,10,11,14)  # A,B,C,D,E,F,G,H,I,J,K,N
foreach ($col in $blankCols21) {
    $cell = $ws.Cells.Item(21, $col)
    $cell.Value = "'"
    $cell.Style = "Normal"
}

# A21 actually carries the card number, like every other row.
$ws.Cells.Item(21, 1).Value = "'12"
$ws.Cells.Item(21, 1).Style = "Normal"

# Date / Serviced by / Correction for the new event.
$ws.Cells.Item(21, 12).Value = "30\9\2024"
$ws.Cells.Item(21, 13).Value = "الخبير"
$ws.Cells.Item(21, 15).Value = "تم تاكيد علي المعيار"

Write-Host "Card12: normalized row 20 and added row 21 event"
